# Scheduled runner update: refresh scraped market-price / profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several Leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 993
$ws.Range("I32").Value = 989.5
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 989.5
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -663.5
$ws.Range("N32").Value = -1652

$ws.Range("H40").Value = 2265.5386
$ws.Range("I40").Value = 1571.4286
$ws.Range("J40").Value = 3075.3333
$ws.Range("K40").Value = 1571.4286
$ws.Range("L40").Value = 3075.3333
$ws.Range("M40").Value = -1396.4286
$ws.Range("N40").Value = -3425.3333

$ws.Range("H76").Value = 3133.3333
$ws.Range("I76").Value = 3141.1765
$ws.Range("J76").Value = 3100
$ws.Range("K76").Value = 3141.1765
$ws.Range("L76").Value = 3100
$ws.Range("M76").Value = -2826.1765
$ws.Range("N76").Value = -3730

$ws.Range("H79").Value = 3133.3333
$ws.Range("I79").Value = 3141.1765
$ws.Range("J79").Value = 3100
$ws.Range("K79").Value = 3141.1765
$ws.Range("L79").Value = 3100
$ws.Range("M79").Value = -2049.1765
$ws.Range("N79").Value = -5284

$ws.Range("H80").Value = 2051.7083
$ws.Range("I80").Value = 1510.4
$ws.Range("J80").Value = 2438.3572
$ws.Range("K80").Value = 4531.200000000001
$ws.Range("L80").Value = 7315.071599999999
$ws.Range("M80").Value = -3533.200000000001
$ws.Range("N80").Value = -9311.071599999999

$ws.Range("H83").Value = 2051.7083
$ws.Range("I83").Value = 1510.4
$ws.Range("J83").Value = 2438.3572
$ws.Range("K83").Value = 13593.6
$ws.Range("L83").Value = 21945.2148
$ws.Range("M83").Value = -8601.6
$ws.Range("N83").Value = -31929.2148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10188.923
$ws.Range("I32").Value = 4610.8716
$ws.Range("K32").Value = 4610.8716
$ws.Range("M32").Value = -4323.8716

$ws.Range("H52").Value = 17999.5
$ws.Range("J52").Value = 17999.5
$ws.Range("L52").Value = 17999.5
$ws.Range("N52").Value = -18635.5

$ws.Range("H63").Value = 3916.2727
$ws.Range("I63").Value = 2520
$ws.Range("J63").Value = 4714.143
$ws.Range("K63").Value = 2520
$ws.Range("L63").Value = 4714.143
$ws.Range("M63").Value = -1834
$ws.Range("N63").Value = -6086.143

$ws.Range("H66").Value = 3916.2727
$ws.Range("I66").Value = 2520
$ws.Range("J66").Value = 4714.143
$ws.Range("K66").Value = 12600
$ws.Range("L66").Value = 23570.715
$ws.Range("M66").Value = -9168
$ws.Range("N66").Value = -30434.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1939.0571
$ws.Range("I105").Value = 1829.862
$ws.Range("J105").Value = 2466.8333
$ws.Range("K105").Value = 1829.862
$ws.Range("L105").Value = 2466.8333
$ws.Range("M105").Value = -82.86200000000008
$ws.Range("N105").Value = -5960.8333

$ws.Range("H107").Value = 2816.6667
$ws.Range("I107").Value = 2780
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2780
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -860
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 578.4167
$ws.Range("I22").Value = 475.25
$ws.Range("J22").Value = 630
$ws.Range("K22").Value = 475.25
$ws.Range("L22").Value = 630
$ws.Range("M22").Value = -125.25
$ws.Range("N22").Value = -1330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44404.086
$ws.Range("I2").Value = 2527.25
$ws.Range("J2").Value = 53220.26
$ws.Range("K2").Value = 15163.5
$ws.Range("L2").Value = 319321.56
$ws.Range("M2").Value = -15050.5
$ws.Range("N2").Value = -319547.56

$ws.Range("H38").Value = 253.81818
$ws.Range("I38").Value = 120
$ws.Range("J38").Value = 283.55554
$ws.Range("K38").Value = 360
$ws.Range("L38").Value = 850.66662
$ws.Range("M38").Value = -13
$ws.Range("N38").Value = -1544.66662

$ws.Range("H97").Value = 630.8333
$ws.Range("J97").Value = 1150
$ws.Range("L97").Value = 3450
$ws.Range("N97").Value = -4442

$ws.Range("H100").Value = 3372.7273
$ws.Range("J100").Value = 3372.7273
$ws.Range("L100").Value = 10118.1819
$ws.Range("N100").Value = -11740.1819

$ws.Range("H109").Value = 2613.2144
$ws.Range("I109").Value = 1725.4
$ws.Range("J109").Value = 3106.4443
$ws.Range("K109").Value = 5176.200000000001
$ws.Range("L109").Value = 9319.332900000001
$ws.Range("M109").Value = -4136.200000000001
$ws.Range("N109").Value = -11399.3329

$ws.Range("H115").Value = 2564.5334
$ws.Range("I115").Value = 1042.6666
$ws.Range("J115").Value = 2945
$ws.Range("K115").Value = 3127.9998
$ws.Range("L115").Value = 8835
$ws.Range("M115").Value = -1952.9998
$ws.Range("N115").Value = -11185

$ws.Range("H140").Value = 2401.0527
$ws.Range("I140").Value = 752.4666999999999
$ws.Range("K140").Value = 2257.4001
$ws.Range("M140").Value = 2922.5999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 21500
$ws.Range("J39").Value = 21500
$ws.Range("L39").Value = 21500
$ws.Range("N39").Value = -22564

$ws.Range("H80").Value = 2233.25
$ws.Range("I80").Value = 2236.2727
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 2236.2727
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -1238.2727
$ws.Range("N80").Value = -4196

$ws.Range("H83").Value = 2233.25
$ws.Range("I83").Value = 2236.2727
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 11181.3635
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -6189.363499999999
$ws.Range("N83").Value = -20984

$ws.Range("H126").Value = 142858240
$ws.Range("I126").Value = 200000860
$ws.Range("J126").Value = 1650
$ws.Range("K126").Value = 600002580
$ws.Range("L126").Value = 4950
$ws.Range("M126").Value = -600000110
$ws.Range("N126").Value = -9890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3099.8572
$ws.Range("I16").Value = 3099.8572
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3099.8572
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2929.8572
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 44286.13
$ws.Range("J22").Value = 889.4737
$ws.Range("L22").Value = 889.4737
$ws.Range("N22").Value = -1479.4737

$ws.Range("H27").Value = 44286.13
$ws.Range("J27").Value = 889.4737
$ws.Range("L27").Value = 889.4737
$ws.Range("N27").Value = -1103.4737

$ws.Range("H46").Value = 969.25
$ws.Range("J46").Value = 827.2308
$ws.Range("L46").Value = 827.2308
$ws.Range("N46").Value = -1203.2308

$ws.Range("H64").Value = 31600
$ws.Range("J64").Value = 31600
$ws.Range("L64").Value = 31600
$ws.Range("N64").Value = -32050

$ws.Range("H67").Value = 31600
$ws.Range("J67").Value = 31600
$ws.Range("L67").Value = 31600
$ws.Range("N67").Value = -33160

$ws.Range("H68").Value = 11900
$ws.Range("I68").Value = 26400
$ws.Range("J68").Value = 2233.3333
$ws.Range("K68").Value = 26400
$ws.Range("L68").Value = 2233.3333
$ws.Range("M68").Value = -25651
$ws.Range("N68").Value = -3731.3333

$ws.Range("H71").Value = 11900
$ws.Range("I71").Value = 26400
$ws.Range("J71").Value = 2233.3333
$ws.Range("K71").Value = 132000
$ws.Range("L71").Value = 11166.6665
$ws.Range("M71").Value = -128256
$ws.Range("N71").Value = -18654.6665

$ws.Range("H82").Value = 1353.091
$ws.Range("J82").Value = 1353.091
$ws.Range("L82").Value = 1353.091
$ws.Range("N82").Value = -2075.091

$ws.Range("H85").Value = 1353.091
$ws.Range("J85").Value = 1353.091
$ws.Range("L85").Value = 1353.091
$ws.Range("N85").Value = -3849.091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1407.5
$ws.Range("I126").Value = 1375
$ws.Range("J126").Value = 1450.8334
$ws.Range("K126").Value = 4125
$ws.Range("L126").Value = 4352.5002
$ws.Range("M126").Value = -1655
$ws.Range("N126").Value = -9292.5002
